$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 = Anke Spijker. Give her 3 keywords: Architectuur (F), Nabewerking (N), Reis (R).
# Copy formatting from existing "marked" cells in the same columns so the new
# cells pick up the same bold / centered / rotated look already used
# throughout the sheet for this kind of mark, then set the value.

$ws.Range("F14").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = "Architectuur"

$ws.Range("N11").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = "Nabewerking"

$ws.Range("R20").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = "Reis"

# Row 8 grew taller to fit the new marks.
$ws.Rows.Item(8).RowHeight = 60

# R7 picks up the same left alignment already used on Q7.
$ws.Range("R7").HorizontalAlignment = $ws.Range("Q7").HorizontalAlignment

# Selection moved.
[void]$ws.Range("AA10").Select()
